$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Merge "   public static void " + "main(" + "String[] " into one run's text.
$find.Execute("public static void main(String[] ", $true, $false, $false, $false, $false, $true, 1, $false, "public static void main(String[] ", 2) | Out-Null

# Merge " = new " + "Scanner( System.in" + " ); " into one run's text.
$find.Execute(" = new Scanner( System.in ); ", $true, $false, $false, $false, $false, $true, 1, $false, " = new Scanner( System.in ); ", 2) | Out-Null
